$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("яблоко")

$ws1.Range("A3").Value = "https://www.apple.com/ru/"
$ws1.Range("B3").Value = "Apple (Россия) – Официальный сайт"
$ws1.Range("C3").Value = "Все инновации Apple, включая iPhone, iPad, Apple Watch, Mac, Apple TV. А также аксессуары, развлечения, справочная информация и многое другое."

$ws1.Range("A4").Value = "https://ru.wikipedia.org/wiki/%D0%AF%D0%B1%D0%BB%D0%BE%D0%BA%D0%BE_(%D0%BF%D0%B0%D1%80%D1%82%D0%B8%D1%8F)"
$ws1.Range("B4").Value = "Яблоко (партия) — Википедия"
$ws1.Range("C4").Value = "Российская объединённая демократическая партия «Яблоко» — зарегистрированная российская политическая партия центристского и ..."

$ws1.Range("A5").Value = "https://ru.wikipedia.org/wiki/%D0%AF%D0%B1%D0%BB%D0%BE%D0%BA%D0%BE"
$ws1.Range("B5").Value = "Яблоко — Википедия"
$ws1.Range("C5").Value = "Я́блоко — сочный плод яблони, который употребляется в пищу в свежем виде, служит сырьём в кулинарии и для приготовления напитков. Наибольшее ..."

$ws1.Range("A6").Value = "https://goldapple.ru/"
$ws1.Range("B6").Value = "«Золотое яблоко» - интернет-магазин косметики и ..."
$ws1.Range("C6").Value = "Косметика и парфюмерия — купите онлайн в интернет-магазине «Золотое яблоко». Более 700 известных брендов: профессиональная, натуральная, ..."

$ws1.Range("A7").Value = "https://vk.com/yabloko_ru"
$ws1.Range("B7").Value = "Партия ЯБЛОКО | ВКонтакте"
$ws1.Range("C7").Value = "Российская объединенная демократическая партия ЯБЛОКО Мы боремся за то, чтобы сделать Россию сильной и современной страной, удобной и ..."

$ws1.Range("A8").Value = "https://mosyabloko.ru/"
$ws1.Range("B8").Value = "Московское ЯБЛОКО: Главная"
$ws1.Range("C8").Value = "В Марьино от партии «Яблоко» на дополнительных выборах в Совет депутатов будут баллотироваться Сергей Запольнов и Данила Столь."

$ws1.Range("A9").Value = "https://lenta.ru/lib/14159780/"
$ws1.Range("B9").Value = "РОДП Яблоко - Lenta.ru"
$ws1.Range("C9").Value = "Имела фракцию в нескольких составах Государственной думы в 1993-2003 годах. Российская объединенная демократическая партия `"Яблоко`" ..."

$ws1.Range("A10").Value = "https://ria.ru/organization_JAbloko/"
$ws1.Range("B10").Value = "Яблоко - последние новости сегодня - РИА Новости"
$ws1.Range("C10").Value = "Яблоко. Читайте последние новости на тему в ленте новостей на сайте РИА Новости. Сенатор Алексей Пушков прокомментировал предложение Яна ..."

$ws1.Range("A11").Value = "https://ru.wikiquote.org/wiki/%D0%AF%D0%B1%D0%BB%D0%BE%D0%BA%D0%BE"
$ws1.Range("B11").Value = "Яблоко — Викицитатник"
$ws1.Range("C11").Value = "Я́блоко — съедобный плод яблони, который употребляют в пищу как в свежем, так и в приготовленном виде: сушёном, мочёном, квашеном, печёном, ..."

$ws2 = $wb.Worksheets.Item("абрикос")

$ws2.Range("A7").Value = "https://foodandmood.com.ua/rid/food/708324-chem-polezen-i-opasen-abrikos"
$ws2.Range("B7").Value = "Чем полезен и опасен абрикос - foodandmood.com.ua"
$ws2.Range("C7").Value = "Выбирайте упругие абрикосы, без повреждений и вмятин. Спелый абрикос очень ароматный и имеет ровную окраску. Избегайте абрикосов с темными ..."

$ws2.Range("A8").Value = "https://foodandmood.com.ua/rid/food/708324-chem-polezen-i-opasen-abrikos"
$ws2.Range("B8").Value = ""
$ws2.Range("C8").Value = ""

$ws2.Range("A9").Value = "https://otvet.mail.ru/question/49528960"
$ws2.Range("B9").Value = ""
$ws2.Range("C9").Value = ""

$ws2.Range("A10").Value = "https://yandex.ru/znatoki/question/garden/chem_otlichaetsia_sliva_ot_abrikosa_82ccaacb/"
$ws2.Range("B10").Value = ""
$ws2.Range("C10").Value = ""

$ws2.Range("A11").Value = "https://forum.derev-grad.ru/plodovie-derevya-f94/skol-ko-let-zhivut-derev-ya-t6442.html"
$ws2.Range("B11").Value = ""
$ws2.Range("C11").Value = ""

